$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells (Modules / Professeurs fields)
$ws.Range("C1").Value = "Chef  Module"
$ws.Range("D1").Value = "Composants"

# Column widths (engine quantizes ColumnWidth in 1/6-char pixel steps;
# these inputs land on the closest achievable stored widths to the
# target 35 and 24.5703125 respectively - 35 is exact, 24.5703125 rounds
# to the nearest reachable value of 24.5)
$ws.Columns.Item(3).ColumnWidth = 34.166666666666664
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668

# Selection
$ws.Range("E8").Select()

# Window size/position
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = 5115
$win.Top = 2760
$win.Width = 15375
$win.Height = 8325
